$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lom3234 = "LOM3234 -  Óptica Física  (Requisito)`n"
$lom3259 = "LOM3259 -  Materiais e Dispositivos Eletrônicos  (Indicação de Conjunto)`n"

# Swap the requisite rows: LOM3259 now appears before LOM3234 (row 24 then row 25)
$ws.Range("B24").Value = $lom3259
$ws.Range("C24").Value = $lom3259
$ws.Range("B25").Value = $lom3234
$ws.Range("C25").Value = $lom3234
